$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 4780.8335
$ws.Range("I107").Value = 4306.364
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 4306.364
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -2386.364
$ws.Range("N107").Value = -13840

$ws.Range("H121").Value = 951.6667
$ws.Range("J121").Value = 912.5
$ws.Range("L121").Value = 2737.5
$ws.Range("N121").Value = -6231.5

$ws.Range("H132").Value = 7411435.5
$ws.Range("J132").Value = 2929.3635
$ws.Range("L132").Value = 8788.0905
$ws.Range("N132").Value = -13848.0905

$ws.Range("H138").Value = 1454.11
$ws.Range("I138").Value = 880.8077
$ws.Range("J138").Value = 1655.5405
$ws.Range("K138").Value = 2642.4231
$ws.Range("L138").Value = 4966.6215
$ws.Range("M138").Value = 2497.5769
$ws.Range("N138").Value = -15246.6215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3390.9875
$ws.Range("I32").Value = 3083.0134
$ws.Range("K32").Value = 3083.0134
$ws.Range("M32").Value = -2796.0134

$ws.Range("H61").Value = 2748.5
$ws.Range("I61").Value = 2798.2
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 2798.2
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -2586.2
$ws.Range("N61").Value = -2924

$ws.Range("H74").Value = 1206.7174
$ws.Range("I74").Value = 705.5333000000001
$ws.Range("J74").Value = 2146.4375
$ws.Range("K74").Value = 705.5333000000001
$ws.Range("L74").Value = 2146.4375
$ws.Range("M74").Value = 168.4666999999999
$ws.Range("N74").Value = -3894.4375

$ws.Range("H77").Value = 1206.7174
$ws.Range("I77").Value = 705.5333000000001
$ws.Range("J77").Value = 2146.4375
$ws.Range("K77").Value = 3527.6665
$ws.Range("L77").Value = 10732.1875
$ws.Range("M77").Value = 840.3334999999997
$ws.Range("N77").Value = -19468.1875

$ws.Range("H102").Value = 16667589
$ws.Range("I102").Value = 20834212
$ws.Range("J102").Value = 1100
$ws.Range("K102").Value = 20834212
$ws.Range("L102").Value = 1100
$ws.Range("M102").Value = -20832590
$ws.Range("N102").Value = -4344

$ws.Range("H110").Value = 1373.25
$ws.Range("I110").Value = 903.7143
$ws.Range("J110").Value = 2468.8333
$ws.Range("K110").Value = 903.7143
$ws.Range("L110").Value = 2468.8333
$ws.Range("M110").Value = 1141.2857
$ws.Range("N110").Value = -6558.8333

$ws.Range("H136").Value = 2748.5
$ws.Range("I136").Value = 2798.2
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 8394.599999999999
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -5844.599999999999
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H86").Value = 4165.9414
$ws.Range("I86").Value = 4301.3125
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 4301.3125
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -3178.3125
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 4165.9414
$ws.Range("I89").Value = 4301.3125
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 21506.5625
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -15890.5625
$ws.Range("N89").Value = -21232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 43479376
$ws.Range("I16").Value = 62500964
$ws.Range("K16").Value = 62500964
$ws.Range("M16").Value = -62500677

$ws.Range("H31").Value = 1331.6666
$ws.Range("I31").Value = 1225.4546
$ws.Range("K31").Value = 1225.4546
$ws.Range("M31").Value = -930.4546

$ws.Range("H34").Value = 1331.6666
$ws.Range("I34").Value = 1225.4546
$ws.Range("K34").Value = 1225.4546
$ws.Range("M34").Value = -1023.4546

$ws.Range("H105").Value = 797.8
$ws.Range("I105").Value = 808
$ws.Range("J105").Value = 769.75
$ws.Range("K105").Value = 808
$ws.Range("L105").Value = 769.75
$ws.Range("M105").Value = 939
$ws.Range("N105").Value = -4263.75

$ws.Range("H113").Value = 43479376
$ws.Range("I113").Value = 62500964
$ws.Range("K113").Value = 62500964
$ws.Range("M113").Value = -62498794

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1444.7931
$ws.Range("I5").Value = 1740.05
$ws.Range("J5").Value = 788.6667
$ws.Range("K5").Value = 5220.15
$ws.Range("L5").Value = 2366.0001
$ws.Range("M5").Value = -5108.15
$ws.Range("N5").Value = -2590.0001

$ws.Range("H97").Value = 175
$ws.Range("I97").Value = 175
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 525
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -29
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 2941.5
$ws.Range("J102").Value = 2941.5
$ws.Range("L102").Value = 8824.5
$ws.Range("N102").Value = -13692.5

$ws.Range("H104").Value = 3712.9443
$ws.Range("I104").Value = 2634
$ws.Range("J104").Value = 4252.4165
$ws.Range("K104").Value = 7902
$ws.Range("L104").Value = 12757.2495
$ws.Range("M104").Value = -5281
$ws.Range("N104").Value = -17999.2495

$ws.Range("H107").Value = 8173.154
$ws.Range("I107").Value = 441
$ws.Range("J107").Value = 11609.667
$ws.Range("K107").Value = 1323
$ws.Range("L107").Value = 34829.001
$ws.Range("M107").Value = 597
$ws.Range("N107").Value = -38669.001

$ws.Range("H112").Value = 13002.077
$ws.Range("I112").Value = 3675.6667
$ws.Range("J112").Value = 15800
$ws.Range("K112").Value = 11027.0001
$ws.Range("L112").Value = 47400
$ws.Range("M112").Value = -9919.000100000001
$ws.Range("N112").Value = -49616

$ws.Range("H135").Value = 1444.7931
$ws.Range("I135").Value = 1740.05
$ws.Range("J135").Value = 788.6667
$ws.Range("K135").Value = 15660.45
$ws.Range("L135").Value = 7098.0003
$ws.Range("M135").Value = -13125.45
$ws.Range("N135").Value = -12168.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 5003000
$ws.Range("J20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("N20").Value = -6490

$ws.Range("H24").Value = 6668666.5
$ws.Range("J24").Value = 6000
$ws.Range("L24").Value = 6000
$ws.Range("N24").Value = -6346

$ws.Range("H107").Value = 704.65216
$ws.Range("I107").Value = 786.0714
$ws.Range("J107").Value = 578
$ws.Range("K107").Value = 786.0714
$ws.Range("L107").Value = 578
$ws.Range("M107").Value = 1133.9286
$ws.Range("N107").Value = -4418

$ws.Range("H132").Value = 2314.7273
$ws.Range("I132").Value = 1350
$ws.Range("J132").Value = 3708.2222
$ws.Range("K132").Value = 4050
$ws.Range("L132").Value = 11124.6666
$ws.Range("M132").Value = -1520
$ws.Range("N132").Value = -16184.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3335116.8
$ws.Range("I20").Value = 10000000
$ws.Range("J20").Value = 2675
$ws.Range("K20").Value = 10000000
$ws.Range("L20").Value = 2675
$ws.Range("M20").Value = -9999774
$ws.Range("N20").Value = -3127

$ws.Range("H68").Value = 2344.875
$ws.Range("J68").Value = 3299.2
$ws.Range("L68").Value = 3299.2
$ws.Range("N68").Value = -4797.2

$ws.Range("H71").Value = 2344.875
$ws.Range("J71").Value = 3299.2
$ws.Range("L71").Value = 16496
$ws.Range("N71").Value = -23984

$ws.Range("H122").Value = 15627363
$ws.Range("I122").Value = 31251750
$ws.Range("K122").Value = 93755250
$ws.Range("M122").Value = -93752800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 20979.75
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 20979.75
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 20979.75
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -21675.75

$ws.Range("H37").Value = 12000
$ws.Range("J37").Value = 12000
$ws.Range("L37").Value = 12000
$ws.Range("N37").Value = -12406

$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = 61
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = 304
$ws.Range("N84").Value = -30608

$ws.Range("H122").Value = 31251762
$ws.Range("I122").Value = 35716016
$ws.Range("K122").Value = 107148048
$ws.Range("M122").Value = -107145598
